$wb = $excel.ActiveWorkbook
$targets = $wb.Worksheets.Item("targets")
$setup = $wb.Worksheets.Item("setup")

$setup.Cells.Item(4, 1).Value = $targets.Cells.Item(1, 1).Value2
$setup.Cells.Item(4, 2).Value = $targets.Cells.Item(1, 2).Value2

$targets.Delete()

$enthalpies = $wb.Worksheets.Item("enthalpies")
$enthalpies.Activate()
$enthalpies.Range("A2").Select()
